# Adds a new "0.7.1" version-history entry (row 25) on Sheet1, re-using the
# "Open points" text from the 0.7.0 row above it, and moves the active
# selection to D26 as in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A25").Value = "0.7.1"
$ws.Range("B25").Value = "AUTOMATA CELULAR - copia (35)"
$ws.Range("C25").Value = $ws.Range("C24").Value2
$ws.Range("D25").Value = "-Reworked Greed."
$ws.Range("E25").Value = $ws.Range("E24").Value2
$ws.Range("F25").Value = $ws.Range("F24").Value2

$ws.Rows.Item(25).RowHeight = 72

$ws.Range("D26").Select()
